$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "SPV"
$ws.Range("B2").Value = "Supervisor"

$ws.Range("B2").Select()
